# Update Price (D) and Volume(1h) (E) columns with latest scraped crypto data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.043.33"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.830.24"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("D4").Value = "'0.9985"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'241.23"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'0.6268"
$ws.Range("E6").Value = "  -4.91%  "
$ws.Range("D8").Value = "'0.07589"
$ws.Range("E8").Value = "  +2.38%  "
$ws.Range("D9").Value = "'45.00"
$ws.Range("E9").Value = "  +7.61%  "
$ws.Range("D10").Value = "'0.2912"
$ws.Range("E10").Value = "  -0.65%  "
$ws.Range("D11").Value = "'22.78"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "'0.07636"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "1.827.86"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "'4.956"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "'0.6652"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "'82.26"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "'0.000009120"
$ws.Range("E17").Value = "  +6.68%  "
$ws.Range("D18").Value = "'5.983"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "28.983.35"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "'224.54"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("D23").Value = "'7.199"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "'159.75"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "'8.412"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "'0.1361"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("D28").Value = "'17.81"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").Value = "'1.495"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "'4.031"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'4.048"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("D33").Value = "'0.05194"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("D34").Value = "'1.844"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").Value = "'0.7321"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("D37").Value = "'2.614"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").Value = "1.279.60"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "'2.762"
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'6.499"
$ws.Range("E41").Value = "  +7.68%  "
$ws.Range("D42").Value = "'0.8922"
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").Value = "'0.9996"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "'101.58"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").Value = "1.977.82"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").Value = "'0.5107"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "'63.85"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'0.3981"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("D50").Value = "'0.07226"
$ws.Range("E50").Value = "  -15.98%  "
$ws.Range("D51").Value = "'8.827"
$ws.Range("E51").Value = "  +1.22%  "
